$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # F2: 1153 -> 1155
    $ws.Range("F2").Value = 1155

    # C4: append （取消） to the event name
    $ws.Range("C4").Value = "合肥·书香璃樱动漫游戏嘉年华（取消）"

    # G4: numeric 50 -> text "不可售"
    $ws.Range("G4").Value = "不可售"

    # F6: 143 -> 145
    $ws.Range("F6").Value = 145

    # F10: 5222 -> 5231
    $ws.Range("F10").Value = 5231

    # F11: 4793 -> 4796
    $ws.Range("F11").Value = 4796

    # F12: 15 -> 16
    $ws.Range("F12").Value = 16
}
